$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{ 1 = "LTQ-Orbitrap_86"; 2 = "LTQ-OrbitrapO_65"; 3 = "LTQ-OrbitrapW_56" }

for ($r = 2; $r -le 46; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = [int]$cell.Value()
    $cell.Value = $map[$val]
}
